$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.674.68"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.699.24"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Formula = "'1.003"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Formula = "'315.49"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Formula = "'1.003"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Formula = "'0.3927"
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("D8").Formula = "'0.4030"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Formula = "'1.002"
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("D11").Formula = "'53.31"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").Formula = "'0.08843"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Formula = "'7.467"
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").Formula = "'23.65"
$ws.Range("E14").Value = "  +1.82%  "
$ws.Range("D15").Formula = "'8.144"
$ws.Range("E15").Value = "  +7.28%  "
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "1.705.30"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Formula = "'99.66"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Formula = "'0.07036"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Formula = "'7.048"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Formula = "'14.69"
$ws.Range("E23").Value = "  +4.57%  "
$ws.Range("D24").Value = "24.675.88"
$ws.Range("D25").Formula = "'3.143"
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("D26").Formula = "'2.364"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").Formula = "'22.65"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("D28").Formula = "'163.02"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").Formula = "'8.729"
$ws.Range("E29").Value = "  +16.68%  "
$ws.Range("D30").Formula = "'135.76"
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Formula = "'5.173"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Formula = "'0.09002"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").Formula = "'7.655"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("D34").Formula = "'1.067"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Formula = "'1.975"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Formula = "'11.05"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Formula = "'0.2751"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Formula = "'14.51"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Formula = "'0.02778"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").Formula = "'0.09120"
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Formula = "'1.461"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Formula = "'0.7672"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Formula = "'15.88"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").Formula = "'0.7165"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Formula = "'2.557"
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("D46").Formula = "'4.216"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Formula = "'1.003"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Formula = "'1.344"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Formula = "'139.68"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Formula = "'0.07975"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").Formula = "'90.19"
$ws.Range("E51").Value = "  +2.47%  "
